$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-0.9253818133688692, 1.084398577801779, 1.0448756794302383, 2.43193904495486)"
$ws.Range("C2").Value = "JSU(-0.864641889850274, 1.049053112163711, 4.560032182629527, 4.920485607784503)"
$ws.Range("D2").Value = "JSU(-0.671348683652208, 0.9856592116870488, 0.9561681775097021, 2.372099718691907)"
$ws.Range("E2").Value = "JSU(-1.027443545856932, 1.0385822068769421, 3.1975052928619627, 4.63010650752927)"
